# 15 Oct 2024 - LV Contacts - Final
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Contact
$ws2 = $wb.Worksheets.Item(2)   # Users
$ws3 = $wb.Worksheets.Item(3)   # Relationship

# ---------------------------------------------------------------
# Sheet "Contact": change the existing External-Contact last name
# and add four new columns (Email / Phone / FullName / ContactType)
# ---------------------------------------------------------------
$ws1.Range("C2").Value = "ExtContact"

$ws1.Range("D1").Value = "Email"
$ws1.Range("E1").Value = "Phone"
$ws1.Range("F1").Value = "FullName"
$ws1.Range("G1").Value = "ContactType"
$ws1.Range("D1:G1").Font.Bold = $true

$ws1.Range("D2").Value = "TestExtContact@email.com"
$ws1.Range("E2").Value = "(541) 754-3010"
$ws1.Range("E2").NumberFormat = "@"
$ws1.Range("F2").Value = "Test ExtContact"
$ws1.Range("G2").Value = "External Contact"

# Hyperlink the new email address cell (also applies the built-in
# "Hyperlink" cell style, like Excel does automatically)
$ws1.Hyperlinks.Add($ws1.Range("D2"), "mailto:TestExtContact@email.com")

# Match the column widths Excel would have auto-fit to the new content
$ws1.Columns.Item(4).ColumnWidth = 22.833333333333336
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.666666666666666
$ws1.Columns.Item(7).ColumnWidth = 13.833333333333332

# ---------------------------------------------------------------
# Sheet "Users": swap the sample user for a new one
# ---------------------------------------------------------------
$ws2.Range("A1").Value = "CF Financial"
$ws2.Range("A2").Value = "Amanda Donovan"

# ---------------------------------------------------------------
# Sheet "Relationship": swap the sample relationship values
# ---------------------------------------------------------------
$ws3.Range("A2").Value = "Houlihan Employee"
$ws3.Range("D2").Value = "Test ExtContact"
$ws3.Columns.Item(1).ColumnWidth = 16

# ---------------------------------------------------------------
# View state: Relationship tab becomes the active / selected tab,
# with C6 selected; Users tab loses the "active" state.
# ---------------------------------------------------------------
$ws2.Range("C6").Select() | Out-Null
$ws3.Activate() | Out-Null
$ws3.Range("C6").Select() | Out-Null

Write-Host "edit complete"
